$d = $word.ActiveDocument

# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
[void]$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                               $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing address paragraph
#       "2968 Lamory Pl, Santa Clara CA 95051"
#    into two separate paragraphs:
#       "2968 Lamory Pl"
#       "Santa Clara, CA 95051"
# First remove the ", Santa Clara CA 95051" tail, leaving just the street
# address in its original paragraph.
$rng = $d.Content
[void]$rng.Find.Execute(", Santa Clara CA 95051")
$rng.Text = ""

# Locate that (now shortened) street-address paragraph...
$addrIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq ("2968 Lamory Pl" + [char]13)) {
        $addrIndex = $i
        break
    }
}

# ...insert a new paragraph right after it (inherits the same paragraph /
# run formatting), and give the new paragraph the city/state/zip text.
$addrPara = $d.Paragraphs($addrIndex)
$addrPara.Range.InsertParagraphAfter()
$cityPara = $d.Paragraphs($addrIndex + 1)
$cityPara.Range.Text = "Santa Clara, CA 95051"

# 3. Remove the blank "No Spacing" paragraph that sits between the
#    "...Board of Directors" signature line and the following Title-styled
#    paragraph.
$bodIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "Board of Directors") {
        $bodIndex = $i
        break
    }
}
$blankPara = $d.Paragraphs($bodIndex + 1)
$blankPara.Range.Delete()
